$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force the Price column (D) to Text format so that numeric-looking
# strings (e.g. "1.000", "307.59") are preserved verbatim as text, matching the
# original inline-string cell content, instead of being auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '27.260.24'
$ws.Range("E2").Value = '  +1.37%  '

# Row 3
$ws.Range("D3").Value = '1.907.16'

# Row 4
$ws.Range("E4").Value = '  +0.09%  '

# Row 5
$ws.Range("D5").Value = '307.59'
$ws.Range("E5").Value = '  +0.54%  '

# Row 6
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  +0.03%  '

# Row 7
$ws.Range("D7").Value = '0.5251'
$ws.Range("E7").Value = '  +2.99%  '

# Row 8
$ws.Range("D8").Value = '0.3787'
$ws.Range("E8").Value = '  +3.54%  '

# Row 9
$ws.Range("D9").Value = '0.07269'
$ws.Range("E9").Value = '  +1.33%  '

# Row 10
$ws.Range("D10").Value = '21.33'
$ws.Range("E10").Value = '  +3.66%  '

# Row 11
$ws.Range("D11").Value = '0.9006'
$ws.Range("E11").Value = '  +1.16%  '

# Row 12
$ws.Range("D12").Value = '0.08133'
$ws.Range("E12").Value = '  +8.63%  '

# Row 13
$ws.Range("D13").Value = '1.912.81'
$ws.Range("E13").Value = '  +2.36%  '

# Row 14
$ws.Range("D14").Value = '95.31'
$ws.Range("E14").Value = '  +0.86%  '

# Row 15
$ws.Range("D15").Value = '5.294'
$ws.Range("E15").Value = '  +1.46%  '

# Row 16
$ws.Range("D16").Value = '1.002'
$ws.Range("E16").Value = '  +0.10%  '

# Row 17
$ws.Range("D17").Value = '0.000008627'
$ws.Range("E17").Value = '  +1.65%  '

# Row 18
$ws.Range("E18").Value = '  +2.53%  '

# Row 19
$ws.Range("D19").Value = '1.001'
$ws.Range("E19").Value = '  +0.07%  '

# Row 20
$ws.Range("D20").Value = '27.315.88'
$ws.Range("E20").Value = '  +1.41%  '

# Row 21
$ws.Range("D21").Value = '5.066'
$ws.Range("E21").Value = '  +1.25%  '

# Row 22
$ws.Range("D22").Value = '2.153.67'
$ws.Range("E22").Value = '  +2.02%  '

# Row 23
$ws.Range("E23").Value = '  +2.92%  '

# Row 24
$ws.Range("D24").Value = '6.459'
$ws.Range("E24").Value = '  +1.40%  '

# Row 25
$ws.Range("D25").Value = '2.310'
$ws.Range("E25").Value = '  +10.93%  '

# Row 26
$ws.Range("D26").Value = '146.32'
$ws.Range("E26").Value = '  -0.86%  '

# Row 27
$ws.Range("D27").Value = '1.748'
$ws.Range("E27").Value = '  -1.76%  '

# Row 28
$ws.Range("E28").Value = '  +1.98%  '

# Row 29
$ws.Range("D29").Value = '115.07'
$ws.Range("E29").Value = '  +1.47%  '

# Row 30
$ws.Range("D30").Value = '4.995'
$ws.Range("E30").Value = '  +6.18%  '

# Row 31
$ws.Range("D31").Value = '4.814'
$ws.Range("E31").Value = '  +2.96%  '

# Row 32
$ws.Range("D32").Value = '0.09232'
$ws.Range("E32").Value = '  +1.18%  '

# Row 33
$ws.Range("D33").Value = '0.8043'
$ws.Range("E33").Value = '  +7.66%  '

# Row 34
$ws.Range("D34").Value = '0.05060'
$ws.Range("E34").Value = '  +0.69%  '

# Row 35
$ws.Range("E35").Value = '  +8.15%  '

# Row 36
$ws.Range("E36").Value = '  +0.78%  '

# Row 37
$ws.Range("D37").Value = '3.330'
$ws.Range("E37").Value = '  +3.38%  '

# Row 38
$ws.Range("D38").Value = '2.577'
$ws.Range("E38").Value = '  +2.68%  '

# Row 39
$ws.Range("D39").Value = '0.5730'
$ws.Range("E39").Value = '  +2.35%  '

# Row 40
$ws.Range("D40").Value = '0.01986'
$ws.Range("E40").Value = '  -0.06%  '

# Row 41
$ws.Range("D41").Value = '1.078'
$ws.Range("E41").Value = '  +0.73%  '

# Row 42
$ws.Range("D42").Value = '119.54'
$ws.Range("E42").Value = '  +3.51%  '

# Row 43
$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D43").Value = '8.978'
$ws.Range("E43").Value = '  +4.80%  '

# Row 44
$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").Value = '6.630'
$ws.Range("E44").Value = '  +0.77%  '

# Row 45
$ws.Range("D45").Value = '0.1515'
$ws.Range("E45").Value = '  +2.17%  '

# Row 46
$ws.Range("D46").Value = '0.4852'
$ws.Range("E46").Value = '  +1.79%  '

# Row 47
$ws.Range("D47").Value = '10.25'
$ws.Range("E47").Value = '  +2.09%  '

# Row 48
$ws.Range("D48").Value = '1.000'
$ws.Range("E48").Value = '  +0.03%  '

# Row 49
$ws.Range("D49").Value = '1.624'
$ws.Range("E49").Value = '  +4.51%  '

# Row 50
$ws.Range("D50").Value = '37.64'
$ws.Range("E50").Value = '  +1.58%  '

# Row 51
$ws.Range("D51").Value = '63.78'
$ws.Range("E51").Value = '  +1.29%  '

# Restore the default ("Normal") cell style on the Price column so the saved
# workbook does not retain an extraneous explicit text-format style definition
# on these cells (they originally had no style override).
$ws.Range("D2:D51").Style = "Normal"
